$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Leetcode row (row 4) problem counts: Easy, Medium, Hard
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 15
$ws.Range("E4").Value = 6

# Move the active selection to E4 (as in the edited file)
$ws.Range("E4").Select()
